# category_renaming.xlsx - "More category work, including Platforms conversion task"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Platforms section ---
# Row 14: Mac -> macOS mapping becomes Mac -> Mac
$ws.Range("B14").Value = "Mac"

# --- Subjects section ---
# Row 27 (Cell Biology): normalize mapped value from "cell biology" to "cells"
$ws.Range("B27").Value = "cells"
# Row 29 (cellular automata): normalize mapped value from "cellular automata" to "cells"
$ws.Range("B29").Value = "cells"

# Row 45 (High School) and Row 50 (Middle School): clear the placeholder
# mapping + question-mark note cells entirely (the <c> elements disappear,
# leaving only the category-name cell in column A)
$ws.Range("B45:C45").ClearContents()
$ws.Range("B50:C50").ClearContents()

# --- Intended ages section ---
# Row 68 (everyone): fix the mapped value's punctuation/wording
$ws.Range("B68").Value = "elementary school, middle school, high school, college"

# --- Selection / view state ---
$ws.Range("B29").Select() | Out-Null
